$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dash example notes to also show the literal "--" / "---" markdown forms
$ws.Range("G2").Value = "&ndash;    or     --"
$ws.Range("G3").Value = "&mdash;    or    ---"

# Reflect the new selection that was active when the workbook was saved
$ws.Activate()
$ws.Range("G2:G3").Select()
